$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# "Available Data" counts for date, total_cases, new_cases, population rows increased from 309 to 310
$ws.Range("B2").Value = 310
$ws.Range("B3").Value = 310
$ws.Range("B4").Value = 310
$ws.Range("B5").Value = 310

# Oporavljeni row: missing data count grew from 78 to 79, pct recalculated
$ws.Range("C6").Value = 79
$ws.Range("D6").Value = 0.341991341991342

# Testirani row: missing data count grew from 78 to 79, pct recalculated
$ws.Range("C7").Value = 79
$ws.Range("D7").Value = 0.341991341991342

# Smrtni sl. row: available data count dropped from 232 to 231, missing data count grew from 77 to 79, pct recalculated
$ws.Range("B8").Value = 231
$ws.Range("C8").Value = 79
$ws.Range("D8").Value = 0.341991341991342
